$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("BP1").Value = "08-sep"

$values = @{
  2 = 0
  3 = 19.380715691229927
  4 = 15.099740854572941
  5 = 17.394037006661399
  6 = 0
  7 = 15.856187346090296
  8 = 12.837896349030995
  9 = 16.580101067648886
  10 = 13.77780448024799
  11 = 11.418459548031638
  12 = 0
  13 = 9.3386824285896406
  14 = 0
  15 = 0
  16 = 11.127152296228767
  17 = 0
  18 = 0
}

foreach ($r in $values.Keys) {
  $ws.Cells.Item($r, 68).Value = $values[$r]
}

$ws.Range("BS6").Select()
